$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.980.02"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "'2.303.83"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'114.21"
$ws.Range("E5").Value = "  +18.33%  "
$ws.Range("D6").Value = "'270.31"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("E10").Value = "  +5.88%  "
$ws.Range("D11").Value = "'0.0952"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "'9.07"
$ws.Range("E12").Value = "  +14.48%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'15.95"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'2.647.78"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "'0.857"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "'2.297.57"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "'43.870.04"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "'6.79"
$ws.Range("E20").Value = "  +9.66%  "
$ws.Range("D21").Value = "'72.53"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").Value = "'233.43"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "'9.71"
$ws.Range("E24").Value = "  +6.19%  "
$ws.Range("D25").Value = "'2.87"
$ws.Range("E25").Value = "  +4.96%  "
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("D28").Value = "'42.18"
$ws.Range("E28").Value = "  +9.01%  "
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "'175.61"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'0.0942"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").Value = "'21.65"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("E37").Value = "  +3.47%  "
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").Value = "'73.90"
$ws.Range("E41").Value = "  +14.97%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'13.81"
$ws.Range("E42").Value = "  +12.56%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").Value = "'6.46"
$ws.Range("E43").Value = "  +24.07%  "
$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D44").Value = "'2.40"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").Value = "'102.75"
$ws.Range("E48").Value = "  +5.32%  "
$ws.Range("D49").Value = "'0.1000"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("D51").Value = "'0.466"
$ws.Range("E51").Value = "  +7.55%  "
